$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-18 Friday" "2024-10-19 Saturday"

Replace-Text "886×4=3544" "540×9=4860"
Replace-Text "595×8=4760" "547×9=4923"
Replace-Text "448×7=3136" "890×7=6230"
Replace-Text "361×6=2166" "791×5=3955"
Replace-Text "980×5=4900" "632×8=5056"

Replace-Text "175×5=875" "838×9=7542"
Replace-Text "304×4=1216" "815×5=4075"
Replace-Text "551×9=4959" "253×2=506"
Replace-Text "705×9=6345" "615×4=2460"
Replace-Text "551×7=3857" "219×3=657"

Replace-Text "718×6=4308" "844×9=7596"
Replace-Text "120×3=360" "437×9=3933"
Replace-Text "405×2=810" "356×8=2848"
Replace-Text "286×7=2002" "447×8=3576"
Replace-Text "692×2=1384" "609×7=4263"

Replace-Text "214×5=1070" "710×5=3550"
Replace-Text "776×4=3104" "257×4=1028"
Replace-Text "504×6=3024" "204×3=612"
Replace-Text "597×3=1791" "110×9=990"
Replace-Text "955×3=2865" "546×3=1638"

Replace-Text "916×6=5496" "508×5=2540"
Replace-Text "775×2=1550" "928×7=6496"
Replace-Text "477×3=1431" "837×3=2511"
Replace-Text "916×2=1832" "716×7=5012"
Replace-Text "225×3=675" "848×3=2544"
